# ResumeAKamlani.docx update
#  1. Nudge the photo up slightly (wp:posOffset -88900 -> -88265 EMU,
#     i.e. Shapes(1).Top -7.00pt -> -6.95pt).
#  2. Merge several runs that were needlessly split (same rPr) back into
#     a single run per paragraph -- this is a pure text/run-structure
#     cleanup, Find/Replace with the concatenated text collapses the
#     run boundaries.
#  3. One of those merges also drops a duplicated "Jquery, " fragment
#     from the Languages line.
#  4. Two new (unused) character styles get minted in styles.xml,
#     mirroring the existing ListLabel7x pairs.

$d = $word.ActiveDocument

# --- 1. floating photo vertical offset -------------------------------
$shp = $d.Shapes(1)
$shp.Top = -6.95

# --- helper -------------------------------------------------------------
function Merge-Text($old, $new) {
    $null = $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                                     $true, 1, $false, $new, 2)
}

# --- 2/3. run merges ----------------------------------------------------
Merge-Text "Technologist within R&D Innovation Group responsible for future Technology and Application Advancement." `
           "Technologist within R&D Innovation Group responsible for future Technology and Application Advancement."

Merge-Text "POC designs focused on Low Cost RTLS Long Range People Tracking" `
           "POC designs focused on Low Cost RTLS Long Range People Tracking"

Merge-Text "Independent Research" `
           "Independent Research"

Merge-Text "Innovator of sports and broadcast media products/solutions." `
           "Innovator of sports and broadcast media products/solutions."

Merge-Text "Embedded Software Consultant for Motorsports (NASCAR Trucks) division." `
           "Embedded Software Consultant for Motorsports (NASCAR Trucks) division."

Merge-Text "Enhance communication signaling between OMAP (OMAP3530) and AVR (Atmel)" `
           "Enhance communication signaling between OMAP (OMAP3530) and AVR (Atmel)"

Merge-Text "From acquisition of Beceem Communications - Cellular division specializing in WiMAX and LTE." `
           "From acquisition of Beceem Communications - Cellular division specializing in WiMAX and LTE."

Merge-Text "Windows Mobile (6.1/6.5/7) BSP software for QSD8650/8250 SnapDragon ~1GHz ARM Cortex (ARMv7-A) architecture processors. Responsible for board bring-up, BSP development, and smartbook reference software design. Collaborate with teams to assess/establish processor requirements, develop new features, and integrate software from previous processor baselines. " `
           "Windows Mobile (6.1/6.5/7) BSP software for QSD8650/8250 SnapDragon ~1GHz ARM Cortex (ARMv7-A) architecture processors. Responsible for board bring-up, BSP development, and smartbook reference software design. Collaborate with teams to assess/establish processor requirements, develop new features, and integrate software from previous processor baselines. "

Merge-Text "Ensured wireless network performance for base stations in the Northeastern region. Utilized parameter thresholds to simulate network performance. Analyzed frequency cell planning, addressed call-processing failure problems, and conducted drive tests. Resolved issues around CDMA, CDPD, and AMPS." `
           "Ensured wireless network performance for base stations in the Northeastern region. Utilized parameter thresholds to simulate network performance. Analyzed frequency cell planning, addressed call-processing failure problems, and conducted drive tests. Resolved issues around CDMA, CDPD, and AMPS."

# this merge also drops the duplicated "Jquery, " text
Merge-Text ": C, C++, Python, R, Javascript, Jquery, Perl, Squirrel, JSON, XML, HTML, CSS, ARM, VHDL" `
           ": C, C++, Python, R, Javascript, Perl, Squirrel, JSON, XML, HTML, CSS, ARM, VHDL"

Merge-Text ": EM BAP RFID (EM4325), BRCM WiMAX/LTE (BCM350, BCM21890), TI WiFi (TINet1100B), Marvell WiFi (88W8381/85), Phillips WiFi (BGW211), CSR Bluecore (Casira), Qualcomm GPSOne " `
           ": EM BAP RFID (EM4325), BRCM WiMAX/LTE (BCM350, BCM21890), TI WiFi (TINet1100B), Marvell WiFi (88W8381/85), Phillips WiFi (BGW211), CSR Bluecore (Casira), Qualcomm GPSOne "

# --- 4. mint the two trailing ListLabel character styles ---------------
$s79 = $d.Styles.Add("ListLabel 79", 2)
$s79.Font.NameBi = "Symbol"

$s80 = $d.Styles.Add("ListLabel 80", 2)
$s80.Font.NameBi = "OpenSymbol"

Write-Output "done"
